$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 534 (pushes old rows 534:570 down to 535:571)
$ws.Rows.Item(534).Insert()

# Populate the new row 534 with the new weekly price record
$ws.Cells.Item(534, 1).Value = 7
$ws.Cells.Item(534, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(534, 3).Value = "Ñuble"
$ws.Cells.Item(534, 4).Value = 45267
$ws.Cells.Item(534, 5).Value = 16
$ws.Cells.Item(534, 6).Value = 100112003
$ws.Cells.Item(534, 7).Value = "Ajo"
$ws.Cells.Item(534, 8).Value = "Chino"
$ws.Cells.Item(534, 9).Value = "Primera"
$ws.Cells.Item(534, 10).Value = 60
$ws.Cells.Item(534, 11).Value = 24000
$ws.Cells.Item(534, 12).Value = 24000
$ws.Cells.Item(534, 13).Value = 24000
$ws.Cells.Item(534, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(534, 15).Value = "China"
$ws.Cells.Item(534, 16).Value = 2400
$ws.Cells.Item(534, 17).Value = 10
$ws.Cells.Item(534, 18).Value = "Hortaliza"
